$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'274.58"
$ws.Range("D3").Value = "'23.02"
$ws.Range("D4").Value = "'6.441"
$ws.Range("D8").Value = "'1.406"
$ws.Range("D9").Value = "'0.8308"
$ws.Range("D11").Value = "'0.1630"
$ws.Range("D12").Value = "'0.08290"
$ws.Range("D13").Value = "'0.03437"
$ws.Range("D14").Value = "'0.03106"
$ws.Range("D15").Value = "'0.09299"
$ws.Range("D16").Value = "'3.875"
$ws.Range("D17").Value = "'0.001649"
$ws.Range("D18").Value = "'0.04789"
$ws.Range("D19").Value = "'0.006407"
$ws.Range("D20").Value = "'0.005676"
$ws.Range("D23").Value = "'3.712"
$ws.Range("D40").Value = "'0.04711"
$ws.Range("D41").Value = "'0.007058"
$ws.Range("D42").Value = "'0.1161"
$ws.Range("D43").Value = "'0.003351"
$ws.Range("D44").Value = "'0.01217"
$ws.Range("D45").Value = "'0.00006271"
$ws.Range("D48").Value = "'0.7967"
$ws.Range("D49").Value = "'0.02871"
$ws.Range("D50").Value = "'0.00002301"
